$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 (A6=5, Shelter): C6 changes from "SHELTER" to the homeowner/renter note
$ws.Range("C6").Value = "910050 (for homeowner), RNTDWELL for renter"

# --- Row 11 (A11=10, New and used vehicles): B11 corrected typo "cares" -> "cars"
$ws.Range("B11").Value = "New and used cars, fees, and maintenance"

# --- New columns F (PCE Line Number) and G (PCE range of lines), plus sparse H notes
$ws.Range("F1").Value = "PCE Line Number"
$ws.Range("G1").Value = "PCE range of lines"

$ws.Range("F2").Value = "229, 72, 101"
$ws.Range("G2").Value = "72-96, 101, 229-242"

$ws.Range("F3").Value = "239, 97"
$ws.Range("G3").Value = "97-100"

$ws.Range("F4").Value = 139
$ws.Range("G4").Value = 139

$ws.Range("F5").Value = "161, 276, 280, 283"
$ws.Range("G5").Value = "161-167, 276-283"
$ws.Range("H5").Value = "Not sure if postal service, line 280 should be here"

$ws.Range("F6").Value = 151
$ws.Range("G6").Value = "151-160, "

$ws.Range("F7").Value = "21 ex 28, 60 ex 61"
$ws.Range("G7").Value = "21-27, 29-35, 60, 67-69"

$ws.Range("F8").Value = 28
$ws.Range("G8").Value = 28

$ws.Range("F9").Value = "102, 61"
$ws.Range("G9").Value = "61-63, 102-110, "

$ws.Range("F10").Value = 195
$ws.Range("G10").Value = "195-204"

$ws.Range("F11").Value = "187, 4"
$ws.Range("G11").Value = "4-20, 187-194"

$ws.Range("F12").Value = "135, 301, 298-300, 309"
$ws.Range("G12").Value = "135-138, 298-320"

$ws.Range("F13").Value = "246, 294"
$ws.Range("G13").Value = "246-273, 294-297"

$ws.Range("F14").Value = "36, 140, 124, 205, 328"
$ws.Range("G14").Value = "36-59, 124-128, 140-142, 205-227, 243-245, 327-335"
$ws.Range("H14").Value = "Make sure to net foreign travel?"

$ws.Range("F15").Value = "129, 321, 293"
$ws.Range("G15").Value = "129-134, 293, 321-326"

$ws.Range("F16").Value = 111
$ws.Range("G16").Value = "111-117"

$ws.Range("F17").Value = "119, 168"
$ws.Range("G17").Value = "64-66, 119-123, 168-185"

$ws.Range("F18").Value = 284
$ws.Range("G18").Value = "284-291"

# --- New rows 20-22 (row 19 intentionally left blank), each a sub-note of
# the NIPA/PCE reconciliation with its own small-caption style
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "      Net expenditures abroad by U.S. residents"
$ws.Range("F20").Value = 143
$ws.Range("G20").Value = "143-147"
$ws.Range("H20").Value = "Need to net these out from non-durable total?"
$r20 = $ws.Range("B20")
$r20.NumberFormat = "#,##0"
$r20.Font.Name = "MS Sans Serif"
$r20.Font.Size = 10

$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "        Less: Expenditures in the United States by nonresidents"
$ws.Range("F21").Value = 332
$ws.Range("G21").Value = "332-335"
$r21 = $ws.Range("B21")
$r21.NumberFormat = "#,##0"
$r21.Font.Name = "MS Sans Serif"
$r21.Font.Size = 10

$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "  Final consumption expenditures of nonprofit institutions serving households (NPISH) "
$ws.Range("F22").Value = 336
$ws.Range("G22").Value = "336-362"
$r22 = $ws.Range("B22")
$r22.NumberFormat = "#,##0"
$r22.Font.Name = "MS Sans Serif"
$r22.Font.Size = 10
$r22.Font.Bold = $true

# --- The footnote moves from row 23 down to row 31 (leaving the block in between empty)
$ws.Range("B23").Value = $null
$ws.Range("B31").Value = "** Note, the most diffuclt grouping of consumption goods to reallocate was Level 2, ""Miscellaneous"""

# --- Selection ends on G14, matching the edited workbook's last active cell
$ws.Range("G14").Select()
